$d = $word.ActiveDocument

# The last paragraph in the document body (before the sectPr) is the one
# ending with "(Tervezési minták alkalmazása.)" -- append a new paragraph
# after it that contains a hyperlink.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertParagraphAfter()

# Range for the newly created (still empty) paragraph.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

$url = "https://www.swiftbysundell.com/articles/sorting-swift-collections/"
$hyperlink = $d.Hyperlinks.Add($newRange, $url, "", "", $url)
$hyperlink.Range.Style = "Hiperhivatkozs"
